$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 31737
$ws.Cells.Item(2, 2).Value = "Pedro Costela"
$ws.Cells.Item(2, 3).Value = "Operações"
$ws.Cells.Item(2, 5).Value = 6
$ws.Cells.Item(2, 6).Value = 45096
$ws.Cells.Item(2, 7).Value = 10474.67

# Row 3
$ws.Cells.Item(3, 1).Value = 60502
$ws.Cells.Item(3, 2).Value = "Ana Clara Rocha"
$ws.Cells.Item(3, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(3, 4).Value = "Consulta médica"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 45102
$ws.Cells.Item(3, 7).Value = 5792.85

# Row 4
$ws.Cells.Item(4, 1).Value = 16611
$ws.Cells.Item(4, 2).Value = "Sofia Teixeira"
$ws.Cells.Item(4, 3).Value = "Recursos Humanos"
$ws.Cells.Item(4, 4).Value = "Doença"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 45078
$ws.Cells.Item(4, 7).Value = 9219.870000000001

# Row 5
$ws.Cells.Item(5, 1).Value = 68983
$ws.Cells.Item(5, 2).Value = "Natália Moura"
$ws.Cells.Item(5, 3).Value = "Financeiro"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 45103
$ws.Cells.Item(5, 7).Value = 7496.83

# Row 6
$ws.Cells.Item(6, 1).Value = 85997
$ws.Cells.Item(6, 2).Value = "Sra. Sabrina Castro"
$ws.Cells.Item(6, 3).Value = "Recursos Humanos"
$ws.Cells.Item(6, 5).Value = 7
$ws.Cells.Item(6, 6).Value = 45091
$ws.Cells.Item(6, 7).Value = 3002.56

# Row 7
$ws.Cells.Item(7, 1).Value = 94002
$ws.Cells.Item(7, 2).Value = "Marina Lima"
$ws.Cells.Item(7, 3).Value = "Vendas"
$ws.Cells.Item(7, 4).Value = "Doença"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 45104
$ws.Cells.Item(7, 7).Value = 12342.06

# Row 8
$ws.Cells.Item(8, 1).Value = 55774
$ws.Cells.Item(8, 2).Value = "Vitória Moreira"
$ws.Cells.Item(8, 4).Value = "Problemas pessoais"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 45081
$ws.Cells.Item(8, 7).Value = 9556.08

# Row 9
$ws.Cells.Item(9, 1).Value = 53011
$ws.Cells.Item(9, 2).Value = "Enzo Ribeiro"
$ws.Cells.Item(9, 3).Value = "TI"
$ws.Cells.Item(9, 4).Value = "Problemas pessoais"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 45099
$ws.Cells.Item(9, 7).Value = 3692.04

# Row 10
$ws.Cells.Item(10, 1).Value = 27172
$ws.Cells.Item(10, 2).Value = "Dr. Pedro Henrique Pereira"
$ws.Cells.Item(10, 3).Value = "Marketing"
$ws.Cells.Item(10, 5).Value = 6
$ws.Cells.Item(10, 6).Value = 45084
$ws.Cells.Item(10, 7).Value = 3850.01

# Row 11
$ws.Cells.Item(11, 1).Value = 59814
$ws.Cells.Item(11, 2).Value = "Sr. Eduardo Barbosa"
$ws.Cells.Item(11, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(11, 4).Value = "Doença"
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = 45085
$ws.Cells.Item(11, 7).Value = 8541.610000000001
